# Update legacy GSC export data:
#  - Drop the oldest date row (2025-10-08) and shift every subsequent
#    date/pages-count row up by one, appending a brand-new trailing
#    date (2026-01-06) with a pages count of 0.0.
#  - This mirrors the effect of the daily export rolling its date window
#    forward by one day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastRow = $ws.UsedRange.Rows.Count   # 91 (header row 1 + 90 data rows)

# Shift column A (Date) and column C (Pages) up by one row for all data
# rows, i.e. row N gets what used to be in row N+1.
for ($r = 2; $r -lt $lastRow; $r++) {
    $nextRow = $r + 1
    $dateCell = $ws.Cells.Item($r, 1)
    # Force the assignment to remain plain text (rather than letting Excel
    # auto-convert the recognizable "yyyy-MM-dd" string into a date serial
    # number), then strip the formatting back off so the cell keeps its
    # original (default) style.
    $dateCell.NumberFormat = "@"
    $dateCell.Value2 = $ws.Cells.Item($nextRow, 1).Value2
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 3).Value2 = $ws.Cells.Item($nextRow, 3).Value2
}

# The final row becomes a brand-new date with a pages count of 0.0.
$lastDateCell = $ws.Cells.Item($lastRow, 1)
$lastDateCell.NumberFormat = "@"
$lastDateCell.Value2 = "2026-01-06"
$lastDateCell.ClearFormats()

$ws.Cells.Item($lastRow, 3).Value2 = 0.0
